# Edit script implementing the diff:
#  - Split "...Jakarta Faces page. To begin, bring down the following projects:"
#    into 3 runs, changing "To begin, bring down" -> "You will be working with"
#  - Replace the three gitlab.com URLs with their bare repo names
#  - Split the "single page app we return null." sentence, adding a comma after "app"
#  - Split the "At this point in the workshop we have seen..." sentence, adding a
#    comma after "workshop"
#  - Add <w:semiHidden/> to the DefaultParagraphFont character style

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "...Jakarta Faces page. To begin, bring down the following projects:"
#    -> "...Jakarta Faces page. " + "You will be working with" + " the following projects:"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("in a POJO called upon from a Jakarta Faces page. To begin, bring down the following projects:")
$full = $d.Range($rng.Start, $rng.End)
$full.Text = "in a POJO called upon from a Jakarta Faces page. "

$r2 = $d.Range($full.End, $full.End)
$r2.InsertAfter("You will be working with")

$r3 = $d.Range($r2.End, $r2.End)
$r3.InsertAfter(" the following projects:")

# ---------------------------------------------------------------------------
# 2) Replace the three gitlab.com project URLs with their bare repo names
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "https://gitlab.com/omniprof/mod_06_servletclient_participant.git",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "mod_06_servletclient_participant", 2) | Out-Null

$d.Content.Find.Execute(
    "https://gitlab.com/omniprof/mod_06_jsfclient_participant.git",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "mod_06_jsfclient_participant", 2) | Out-Null

$d.Content.Find.Execute(
    "https://gitlab.com/omniprof/mod_06_restserver.git",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "mod_06_restserver", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "...single page app we return null. " -> "...single page " + "app," + " we return null. "
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute(" The only difference is that JSF can access beans and call action methods directly from the page. Action methods must return a string or null. The string could be the next page to navigate to but as this is a single page app we return null. ")
$full2 = $d.Range($rng2.Start, $rng2.End)
$full2.Text = " The only difference is that JSF can access beans and call action methods directly from the page. Action methods must return a string or null. The string could be the next page to navigate to but as this is a single page "

$r2b = $d.Range($full2.End, $full2.End)
$r2b.InsertAfter("app,")

$r3b = $d.Range($r2b.End, $r2b.End)
$r3b.InsertAfter(" we return null. ")

# ---------------------------------------------------------------------------
# 4) "At this point in the workshop we have seen..." -> "At this point in the " +
#    "workshop," + " we have seen..."
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("At this point in the workshop we have seen how we can construct REST services that can run on a server or in a desktop app. We have seen desktop, Servlet, and JSF clients. Coming up we will learn how to construct a service that accepts a binary file such as an image.")
$full3 = $d.Range($rng3.Start, $rng3.End)
$full3.Text = "At this point in the "

$r2c = $d.Range($full3.End, $full3.End)
$r2c.InsertAfter("workshop,")

$r3c = $d.Range($r2c.End, $r2c.End)
$r3c.InsertAfter(" we have seen how we can construct REST services that can run on a server or in a desktop app. We have seen desktop, Servlet, and JSF clients. Coming up we will learn how to construct a service that accepts a binary file such as an image.")

# ---------------------------------------------------------------------------
# 5) Add <w:semiHidden/> to the DefaultParagraphFont character style
# ---------------------------------------------------------------------------
$styles = $d.Styles
$dpf = $styles.Item("Default Paragraph Font")
$dpf.SemiHidden = $true

Write-Output "done"
